$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column F header (row 1) ---
$ws.Range("F1").Value = "Thời gian seminar"

# --- F2 / F3: seminar dates, bold Times New Roman, text format, bordered, left/center aligned ---
$dates = $ws.Range("F2:F3")
$dates.NumberFormat = "@"
$dates.Font.Name = "Times New Roman"
$dates.Font.Size = 11
$dates.Font.Bold = $true
$dates.HorizontalAlignment = -4131
$dates.VerticalAlignment = -4108
$dates.Borders.LineStyle = 1

$ws.Range("F2").Value = "28/04/2017"
$ws.Range("F3").Value = "05/05/2017"

# --- F4 / F5: blank cells pre-formatted for future entries (same text format, not bold) ---
$blanks = $ws.Range("F4:F5")
$blanks.NumberFormat = "@"
$blanks.Font.Name = "Times New Roman"
$blanks.Font.Size = 11
$blanks.HorizontalAlignment = -4131
$blanks.VerticalAlignment = -4108
$blanks.Borders.LineStyle = 1

# --- Column F width ---
$ws.Columns.Item(6).ColumnWidth = 15.3

# --- Row 3 height shrinks once the new column layout is in place ---
$ws.Rows.Item(3).RowHeight = 30

# --- Final selection, as left by the editing session ---
[void]$ws.Range("F9").Select()
